# Apply the recorded edit:
#  - Fix the "MODEL_CONDITION" header text -> "MODELCONDITION"
#  - Shift the whole table one column to the left (old B:F -> new A:E),
#    which drops the now-empty column A/F and naturally carries each
#    cell's original value/type/style with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text in place first (value-only edit, keeps its style).
$ws.Cells.Item(1, 5).Value2 = "MODELCONDITION"

# Delete column A (rows 1:3) and shift the remaining cells (old B:F) left.
$ws.Range("A1:A3").Delete(-4159)
